# Add DemoBlaze Project Folder
# ------------------------------------------------------------
# The screenshots referenced from the "Bug report" sheet used to live in an
# absolute local folder (C:\Users\ADMIN\Pictures\Screenshots\...). The
# project now ships its own "DemoBlaze_V-1.1\screenshots\" folder, so the
# displayed screenshot paths are updated to point there instead.
# ------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$bugReport = $wb.Worksheets.Item("Bug report")

$bugReport.Range("M3").Value  = "DemoBlaze_V-1.1\screenshots\Sign-Up.png"
$bugReport.Range("M4").Value  = "DemoBlaze_V-1.1\screenshots\Sign-Up.png"
$bugReport.Range("M5").Value  = "DemoBlaze_V-1.1\screenshots\Login.png"
$bugReport.Range("M6").Value  = "DemoBlaze_V-1.1\screenshots\LogoutwithBrowserbackbtn.png"
$bugReport.Range("M7").Value  = "DemoBlaze_V-1.1\screenshots\Unsuccessful Login Attempts.png"
$bugReport.Range("M8").Value  = "DemoBlaze_V-1.1\screenshots\Login.png"
$bugReport.Range("M9").Value  = "DemoBlaze_V-1.1\screenshots\Contact.png"
$bugReport.Range("M11").Value = "DemoBlaze_V-1.1\screenshots\Contact with Invalid Email.png"
$bugReport.Range("M12").Value = "DemoBlaze_V-1.1\screenshots\Contact Message limit.png"
$bugReport.Range("M13").Value = "DemoBlaze_V-1.1\screenshots\cart in quantity and total price.png"
$bugReport.Range("M14").Value = "DemoBlaze_V-1.1\screenshots\Remove Item from Cart.png"
$bugReport.Range("M15").Value = "DemoBlaze_V-1.1\screenshots\cart in quantity and total price.png"
$bugReport.Range("M16").Value = "DemoBlaze_V-1.1\screenshots\Phones Category.png"
$bugReport.Range("M17").Value = "DemoBlaze_V-1.1\screenshots\Laptop Category.png"

# The Contact sheet's message-limit row no longer needs quite as much
# vertical space once re-measured, so its height shrinks slightly.
$contact = $wb.Worksheets.Item("Contact")
$contact.Rows.Item(6).RowHeight = 357

# Move the window/focus onto the "Bug report" sheet (it becomes the active
# tab), with the Screenshots column in view and M17's replacement cell
# (now P17, after a couple of extra columns were added while documenting
# the new screenshots) selected.
$win = $wb.Windows.Item(1)
$win.Left = -110
$win.Top = -110
$win.Width = 19420
$win.Height = 10300

$bugReport.Activate()
$bugReport.Range("P17").Select()
